# Update cryptos list values (Price column D, Volume(1h) column E)
# generated from commit "Updated cryptos list on Wed Nov  8 04:29:10 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price (D) column while keeping it stored as
# plain text (these price strings use "." as a thousands separator and would
# otherwise be auto-converted to numbers by Excel, which would also silently
# truncate trailing zeros, e.g. "8.50" -> 8.5).
function Set-PriceText($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-PriceText $ws.Range("D2") "35.307.68"
$ws.Range("E2").Value = "  +0.50%  "
Set-PriceText $ws.Range("D3") "1.881.51"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.65%  "
Set-PriceText $ws.Range("D5") "245.13"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("E6").Value = "  -0.71%  "
Set-PriceText $ws.Range("D8") "43.49"
$ws.Range("E8").Value = "  +5.47%  "
Set-PriceText $ws.Range("D9") "0.355"
$ws.Range("E9").Value = "  -0.97%  "
Set-PriceText $ws.Range("D10") "53.33"
$ws.Range("E10").Value = "  +0.99%  "
Set-PriceText $ws.Range("D11") "0.0740"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  -1.00%  "
Set-PriceText $ws.Range("D13") "13.33"
$ws.Range("E13").Value = "  +2.17%  "
Set-PriceText $ws.Range("D14") "2.155.30"
$ws.Range("E14").Value = "  -1.05%  "
Set-PriceText $ws.Range("D15") "0.756"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("E16").Value = "  -1.68%  "
Set-PriceText $ws.Range("D17") "1.893.44"
$ws.Range("E17").Value = "  -0.40%  "
Set-PriceText $ws.Range("D18") "35.411.67"
$ws.Range("E18").Value = "  +0.80%  "
Set-PriceText $ws.Range("D19") "72.89"
$ws.Range("E20").Value = "  -1.41%  "
Set-PriceText $ws.Range("D21") "244.24"
$ws.Range("E21").Value = "  +0.57%  "
Set-PriceText $ws.Range("D22") "12.78"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -2.29%  "
Set-PriceText $ws.Range("D24") "2.65"
$ws.Range("E24").Value = "  +8.79%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -6.66%  "
Set-PriceText $ws.Range("D27") "165.36"
$ws.Range("E27").Value = "  -0.98%  "
Set-PriceText $ws.Range("D28") "8.50"
$ws.Range("E28").Value = "  -0.55%  "
Set-PriceText $ws.Range("D29") "18.28"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  -1.97%  "
Set-PriceText $ws.Range("D31") "4.128.46"
$ws.Range("E31").Value = "  +0.01%  "
Set-PriceText $ws.Range("D32") "1.69"
$ws.Range("E32").Value = "  +7.58%  "
Set-PriceText $ws.Range("D33") "4.26"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  -7.08%  "
Set-PriceText $ws.Range("D36") "4.14"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").Value = "  -0.72%  "
Set-PriceText $ws.Range("D38") "0.842"
$ws.Range("E38").Value = "  -1.33%  "
Set-PriceText $ws.Range("D40") "0.0696"
$ws.Range("E40").Value = "  +7.08%  "
Set-PriceText $ws.Range("D41") "17.18"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("E42").Value = "  +1.19%  "
Set-PriceText $ws.Range("D43") "96.07"
$ws.Range("E43").Value = "  -6.88%  "
$ws.Range("E44").Value = "  -2.48%  "
Set-PriceText $ws.Range("D45") "1.303.18"
$ws.Range("E45").Value = "  -1.24%  "
Set-PriceText $ws.Range("D46") "2.32"
$ws.Range("E46").Value = "  -4.54%  "
Set-PriceText $ws.Range("D47") "0.0795"
$ws.Range("E47").Value = "  +6.63%  "
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("E49").Value = "  -0.84%  "
Set-PriceText $ws.Range("D50") "12.14"
$ws.Range("E50").Value = "  +1.59%  "
Set-PriceText $ws.Range("D51") "6.23"
$ws.Range("E51").Value = "  -5.48%  "
